$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.37"
$ws.Range("D4").Value = "'5.256"
$ws.Range("D5").Value = "'0.05822"
$ws.Range("D6").Value = "'6.457"
$ws.Range("D8").Value = "'0.8082"
$ws.Range("D9").Value = "'0.9005"
$ws.Range("D10").Value = "'0.1376"
$ws.Range("D11").Value = "'0.07072"
$ws.Range("D12").Value = "'0.03081"
$ws.Range("D13").Value = "'0.03031"
$ws.Range("D14").Value = "'0.09314"
$ws.Range("D15").Value = "'3.818"
$ws.Range("D16").Value = "'0.001541"
$ws.Range("D17").Value = "'0.04709"
$ws.Range("D18").Value = "'0.0006010"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006222"
$ws.Range("D20").Value = "'0.001261"
$ws.Range("D21").Value = "'0.004056"
$ws.Range("D22").Value = "'0.00008702"
$ws.Range("D23").Value = "'3.559"
$ws.Range("D24").Value = "'2.178"
$ws.Range("D25").Value = "'0.3169"
$ws.Range("D26").Value = "'0.1318"
$ws.Range("D40").Value = "'0.03808"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1051"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002522"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003239"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("D44").Value = "'0.006948"
$ws.Range("D45").Value = "'0.00005315"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.5110"
$ws.Range("D48").Value = "'0.007015"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
